$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the "Acelga / Vega Central Mapocho" block
# (this pushes the existing rows 980-1046 down to 982-1048, matching the
# new sheet dimension A1:R1048).
$ws.Range("A980:A981").EntireRow.Insert()

# New row 980 - weekly price record (Primera)
$ws.Range("A980").Value = 9
$ws.Range("B980").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C980").Value = "Metropolitana"
$ws.Range("D980").Value = 45265
$ws.Range("E980").Value = 13
$ws.Range("F980").Value = 100112009
$ws.Range("G980").Value = "Acelga"
$ws.Range("H980").Value = "Sin especificar"
$ws.Range("I980").Value = "Primera"
$ws.Range("J980").Value = 70
$ws.Range("K980").Value = 20000
$ws.Range("L980").Value = 22000
$ws.Range("M980").Value = 21000
$ws.Range("N980").Value = "$/docena de atados"
$ws.Range("O980").Value = "Región Metropolitana"
$ws.Range("P980").Value = 7000
$ws.Range("Q980").Value = 3
$ws.Range("R980").Value = "Hortaliza"

# New row 981 - weekly price record (Segunda)
$ws.Range("A981").Value = 9
$ws.Range("B981").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C981").Value = "Metropolitana"
$ws.Range("D981").Value = 45265
$ws.Range("E981").Value = 13
$ws.Range("F981").Value = 100112009
$ws.Range("G981").Value = "Acelga"
$ws.Range("H981").Value = "Sin especificar"
$ws.Range("I981").Value = "Segunda"
$ws.Range("J981").Value = 52
$ws.Range("K981").Value = 18000
$ws.Range("L981").Value = 18000
$ws.Range("M981").Value = 18000
$ws.Range("N981").Value = "$/docena de atados"
$ws.Range("O981").Value = "Región Metropolitana"
$ws.Range("P981").Value = 6000
$ws.Range("Q981").Value = 3
$ws.Range("R981").Value = "Hortaliza"
